$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.06842864744397358
$ws.Range("J2").Value = 0.06842864744397358
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.07259900000000001
$ws.Range("N2").Value = 0.217797
$ws.Range("O2").Value = 0.0162094769588191
$ws.Range("P2").Value = 0.0162094769588191
$ws.Range("Q2").Value = 0.005551911726333335
$ws.Range("R2").Value = 0.04996720553700001
$ws.Range("S2").Value = 0.001109192584066245
$ws.Range("T2").Value = 0.001109192584066245

# Row 3
$ws.Range("I3").Value = 0.06842864744397358
$ws.Range("J3").Value = 0.06842864744397358
$ws.Range("O3").Value = 0.9349228167457665
$ws.Range("P3").Value = 0.9349228167457664
$ws.Range("S3").Value = 0.06397550381442277
$ws.Range("T3").Value = 0.06397550381442277

# Row 4
$ws.Range("I4").Value = 0.06842864744397358
$ws.Range("J4").Value = 0.06842864744397358
$ws.Range("M4").Value = 0.2188686666666667
$ws.Range("N4").Value = 0.656606
$ws.Range("O4").Value = 0.04886770629541442
$ws.Range("P4").Value = 0.04886770629541441
$ws.Range("Q4").Value = 0.01673768945844445
$ws.Range("R4").Value = 0.150639205126
$ws.Range("S4").Value = 0.003343951045484561
$ws.Range("T4").Value = 0.003343951045484561

# Row 5
$ws.Range("G5").Value = 0.621785
$ws.Range("H5").Value = 1.865355
$ws.Range("I5").Value = 0.5563733034589394
$ws.Range("J5").Value = 0.5563733034589394
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.07259900000000001
$ws.Range("N5").Value = 0.217797
$ws.Range("O5").Value = 0.0162094769588191
$ws.Range("P5").Value = 0.0162094769588191
$ws.Range("Q5").Value = 0.04514096921500001
$ws.Range("R5").Value = 0.406268722935
$ws.Range("S5").Value = 0.009018520242919745
$ws.Range("T5").Value = 0.009018520242919745

# Row 6
$ws.Range("G6").Value = 0.621785
$ws.Range("H6").Value = 1.865355
$ws.Range("I6").Value = 0.5563733034589394
$ws.Range("J6").Value = 0.5563733034589394
$ws.Range("O6").Value = 0.9349228167457665
$ws.Range("P6").Value = 0.9349228167457664
$ws.Range("Q6").Value = 2.60362022762
$ws.Range("R6").Value = 23.43258204858
$ws.Range("S6").Value = 0.5201660960319787
$ws.Range("T6").Value = 0.5201660960319787

# Row 7
$ws.Range("G7").Value = 0.621785
$ws.Range("H7").Value = 1.865355
$ws.Range("I7").Value = 0.5563733034589394
$ws.Range("J7").Value = 0.5563733034589394
$ws.Range("M7").Value = 0.2188686666666667
$ws.Range("N7").Value = 0.656606
$ws.Range("O7").Value = 0.04886770629541442
$ws.Range("P7").Value = 0.04886770629541441
$ws.Range("Q7").Value = 0.1360892539033333
$ws.Range("R7").Value = 1.22480328513
$ws.Range("S7").Value = 0.02718868718404093
$ws.Range("T7").Value = 0.02718868718404093

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.4193093333333334
$ws.Range("H8").Value = 1.257928
$ws.Range("I8").Value = 0.375198049097087
$ws.Range("J8").Value = 0.375198049097087
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07259900000000001
$ws.Range("N8").Value = 0.217797
$ws.Range("O8").Value = 0.0162094769588191
$ws.Range("P8").Value = 0.0162094769588191
$ws.Range("Q8").Value = 0.03044143829066667
$ws.Range("R8").Value = 0.2739729446160001
$ws.Range("S8").Value = 0.006081764131833109
$ws.Range("T8").Value = 0.00608176413183311

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.4193093333333334
$ws.Range("H9").Value = 1.257928
$ws.Range("I9").Value = 0.375198049097087
$ws.Range("J9").Value = 0.375198049097087
$ws.Range("O9").Value = 0.9349228167457665
$ws.Range("P9").Value = 0.9349228167457664
$ws.Range("Q9").Value = 1.755787389365334
$ws.Range("R9").Value = 15.802086504288
$ws.Range("S9").Value = 0.350781216899365
$ws.Range("T9").Value = 0.350781216899365

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.4193093333333334
$ws.Range("H10").Value = 1.257928
$ws.Range("I10").Value = 0.375198049097087
$ws.Range("J10").Value = 0.375198049097087
$ws.Range("M10").Value = 0.2188686666666667
$ws.Range("N10").Value = 0.656606
$ws.Range("O10").Value = 0.04886770629541442
$ws.Range("P10").Value = 0.04886770629541441
$ws.Range("Q10").Value = 0.09177367470755557
$ws.Range("R10").Value = 0.8259630723680002
$ws.Range("S10").Value = 0.01833506806588893
$ws.Range("T10").Value = 0.01833506806588893
